$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44307
$ws.Range("M2").Value = 40

# Row 5
$ws.Range("D5").Value = 44314
$ws.Range("M5").Value = 47
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 9000
$ws.Range("S5").Value = 900

# Row 6
$ws.Range("D6").Value = 44326
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 1000

# Row 7
$ws.Range("D7").Value = 44326
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 67
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("S7").Value = 800

# Row 8
$ws.Range("D8").Value = 44321
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 58
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 9000
$ws.Range("P8").Value = 9000
$ws.Range("S8").Value = 900

# Row 9
$ws.Range("D9").Value = 44315
$ws.Range("M9").Value = 45

# Row 10
$ws.Range("D10").Value = 44333
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 58

# Row 11
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("S11").Value = 900

# Row 12
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("S12").Value = 800

# Row 13
$ws.Range("D13").Value = 44302
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 1000

# Row 14
$ws.Range("D14").Value = 44309
$ws.Range("M14").Value = 45

# Row 15
$ws.Range("D15").Value = 44306
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("S15").Value = 1000

# Row 18
$ws.Range("D18").Value = 44322
$ws.Range("M18").Value = 56

# Row 19
$ws.Range("D19").Value = 44322
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("S19").Value = 800

# Row 20
$ws.Range("D20").Value = 44308

# Row 21
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 48
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 800

# Row 22
$ws.Range("D22").Value = 44301
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 45
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("S22").Value = 1000

# Row 25
$ws.Range("D25").Value = 44312
$ws.Range("M25").Value = 48

# Row 26
$ws.Range("D26").Value = 44319
$ws.Range("M26").Value = 68
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("S26").Value = 1000

# Row 27
$ws.Range("D27").Value = 44319
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 57
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 8000
$ws.Range("R27").Value = "Provincia de Quillota"
$ws.Range("S27").Value = 800

# Row 28
$ws.Range("L28").Value = "Especial"
$ws.Range("M28").Value = 47
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("S28").Value = 1000

# Row 29
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 9000
$ws.Range("O29").Value = 9000
$ws.Range("P29").Value = 9000
$ws.Range("S29").Value = 900

# Row 30
$ws.Range("D30").Value = 44343
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 58
$ws.Range("R30").Value = "Región Metropolitana"

# Row 31
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 45
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("S31").Value = 800

# Row 32
$ws.Range("D32").Value = 44328
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 48
$ws.Range("N32").Value = 7000
$ws.Range("O32").Value = 7000
$ws.Range("P32").Value = 7000
